$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.914.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.60%  "

$ws.Range("D3").Value = "'1.887.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.30%  "

$ws.Range("E4").Value = "  -1.29%  "

$ws.Range("D5").Value = "'326.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.78%  "

$ws.Range("E6").Value = "  -1.05%  "

$ws.Range("D7").Value = "'0.4578"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.90%  "

$ws.Range("D8").Value = "'0.3919"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.73%  "

$ws.Range("D9").Value = "'49.20"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -8.41%  "

$ws.Range("D10").Value = "'0.08210"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.32%  "

$ws.Range("D11").Value = "'1.034"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.43%  "

$ws.Range("D12").Value = "'21.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.00%  "

$ws.Range("D13").Value = "'1.911.10"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.37%  "

$ws.Range("D14").Value = "'7.294"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.21%  "

$ws.Range("D15").Value = "'5.961"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.90%  "

$ws.Range("E16").Value = "  -1.14%  "

$ws.Range("D17").Value = "'88.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.56%  "

$ws.Range("D18").Value = "'0.00001054"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.87%  "

$ws.Range("D19").Value = "'0.06572"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.75%  "

$ws.Range("D20").Value = "'17.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.69%  "

$ws.Range("E21").Value = "  -1.17%  "

$ws.Range("D22").Value = "'5.627"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.10%  "

$ws.Range("D23").Value = "'27.923.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.66%  "

$ws.Range("D24").Value = "'11.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.11%  "

$ws.Range("D25").Value = "'2.305"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.71%  "

$ws.Range("D26").Value = "'2.118.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.35%  "

$ws.Range("D27").Value = "'154.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.45%  "

$ws.Range("D28").Value = "'19.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.59%  "

$ws.Range("D29").Value = "'5.687"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.24%  "

$ws.Range("D30").Value = "'2.099"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.41%  "

$ws.Range("D31").Value = "'123.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.59%  "

$ws.Range("D32").Value = "'0.09520"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.65%  "

$ws.Range("D33").Value = "'0.9535"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.74%  "

$ws.Range("D34").Value = "'1.470"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.67%  "

$ws.Range("D35").Value = "'3.633"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.98%  "

$ws.Range("D36").Value = "'5.444"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.63%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02278"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.44%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.249"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.27%  "

$ws.Range("D39").Value = "'0.06088"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.08%  "

$ws.Range("D40").Value = "'8.546"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.96%  "

$ws.Range("D41").Value = "'0.6093"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.93%  "

$ws.Range("D42").Value = "'1.002"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.07%  "

$ws.Range("D43").Value = "'10.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.48%  "

$ws.Range("D44").Value = "'0.1891"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.43%  "

$ws.Range("D45").Value = "'1.302"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.16%  "

$ws.Range("D46").Value = "'0.5800"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.11%  "

$ws.Range("D47").Value = "'12.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.58%  "

$ws.Range("D48").Value = "'1.983"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.66%  "

$ws.Range("D49").Value = "'3.420"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.19%  "

$ws.Range("D50").Value = "'0.06900"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.12%  "

$ws.Range("D51").Value = "'110.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.77%  "
